# Fruta / hortaliza, semanal
# Reorder the weekly price rows (2, 4, 6) by swapping their date/price data
# cyclically: row2 <- row6, row4 <- row2(old), row6 <- row4(old).
# Rows 3 and 5 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for rows 2, 4 and 6 (the ones that move).
$row2 = @{
    D = $ws.Range("D2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    P = $ws.Range("P2").Value2
}
$row4 = @{
    D = $ws.Range("D4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}
$row6 = @{
    D = $ws.Range("D6").Value2
    J = $ws.Range("J6").Value2
    K = $ws.Range("K6").Value2
    L = $ws.Range("L6").Value2
    M = $ws.Range("M6").Value2
    P = $ws.Range("P6").Value2
}

# Row 2 gets what used to be in row 6.
$ws.Range("D2").Value = $row6.D
$ws.Range("J2").Value = $row6.J
$ws.Range("K2").Value = $row6.K
$ws.Range("L2").Value = $row6.L
$ws.Range("M2").Value = $row6.M
$ws.Range("P2").Value = $row6.P

# Row 4 gets what used to be in row 2.
$ws.Range("D4").Value = $row2.D
$ws.Range("J4").Value = $row2.J
$ws.Range("K4").Value = $row2.K
$ws.Range("L4").Value = $row2.L
$ws.Range("M4").Value = $row2.M
$ws.Range("P4").Value = $row2.P

# Row 6 gets what used to be in row 4.
$ws.Range("D6").Value = $row4.D
$ws.Range("J6").Value = $row4.J
$ws.Range("K6").Value = $row4.K
$ws.Range("L6").Value = $row4.L
$ws.Range("M6").Value = $row4.M
$ws.Range("P6").Value = $row4.P
